# Insert a new weekly record for "Fruta Vega Modelo de Temuco - Coco":
# a new row is inserted at row 48, pushing the existing rows 48:71 down
# to 49:72. The new row re-uses the constant descriptive columns
# (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) shared by every record in this sheet and
# carries its own Fecha/Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 48:71 down to 49:72, inserting a blank row 48.
$ws.Rows.Item(48).Insert()

# Populate the new row 48 with the new record. The descriptive columns are
# identical for every row of this product (Coco, Vega Modelo de Temuco),
# so copy them straight from the row right below (the former row 48,
# now shifted to row 49).
$ws.Range("A48").Value = $ws.Range("A49").Value()
$ws.Range("B48").Value = $ws.Range("B49").Value()
$ws.Range("C48").Value = $ws.Range("C49").Value()
$ws.Range("D48").Value = 44726
$ws.Range("E48").Value = $ws.Range("E49").Value()
$ws.Range("F48").Value = $ws.Range("F49").Value()
$ws.Range("G48").Value = $ws.Range("G49").Value()
$ws.Range("H48").Value = $ws.Range("H49").Value()
$ws.Range("I48").Value = $ws.Range("I49").Value()
$ws.Range("J48").Value = $ws.Range("J49").Value()
$ws.Range("K48").Value = $ws.Range("K49").Value()
$ws.Range("L48").Value = $ws.Range("L49").Value()
$ws.Range("M48").Value = 20
$ws.Range("N48").Value = 30000
$ws.Range("O48").Value = 30000
$ws.Range("P48").Value = 30000
$ws.Range("Q48").Value = $ws.Range("Q49").Value()
$ws.Range("R48").Value = $ws.Range("R49").Value()
$ws.Range("S48").Value = 1500
$ws.Range("T48").Value = $ws.Range("T49").Value()

# Match the date-column number format used by every other row.
$ws.Range("D48").NumberFormat = $ws.Range("D49").NumberFormat()
